$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title
#    paragraph. We seed it from the existing bold "Play Biergarten Fest
#    Free..." paragraph (near the end of the doc) because that paragraph
#    already has the exact <w:r/><w:r><w:rPr><w:b/>...> run shape we need -
#    copying its FormattedText preserves the empty leading run, whereas
#    typing straight into Range.Text / Find-Replace collapses it away.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$boldSourcePara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$boldSourceText = "Play Biergarten Fest Free: Lively Slot Game Overview"
$boldSource = $d.Range($boldSourcePara.Range.Start, $boldSourcePara.Range.End - 1)
$metaPara.Range.FormattedText = $boldSource.FormattedText

# Overwrite "Play Biergarten Fest Free: Lively Slot Game Overview" -> "Meta description"
# using a directly-addressed Range (not Find/Replace) so the leading empty
# run stays untouched.
$boldRunStart = $metaPara.Range.Start
$boldRunEnd = $boldRunStart + $boldSourceText.Length
$boldRun = $d.Range($boldRunStart, $boldRunEnd)
$boldRun.Text = "Meta description"

# Seed a plain (no rPr) trailing run by copying formatting from an existing
# un-styled run elsewhere in the doc, then overwrite its text - same trick
# as above, keeps run boundaries intact instead of merging into the bold run.
$plainSourcePara = $d.Paragraphs.Item(5)
$plainSeedLen = 5
$plainSource = $d.Range($plainSourcePara.Range.Start, $plainSourcePara.Range.Start + $plainSeedLen)

$tailInsertPos = $metaPara.Range.End - 1
$tailRange = $d.Range($tailInsertPos, $tailInsertPos)
$tailRange.FormattedText = $plainSource.FormattedText

$metaDescText = ": Discover Biergarten Fest slot game, inspired by Oktoberfest. Enjoy Wilds, Scatters, and bonus round. Play for free, and experience immersive graphics."
$tailRun = $d.Range($tailInsertPos, $tailInsertPos + $plainSeedLen)
$tailRun.Text = $metaDescText

# ---------------------------------------------------------------------------
# 2) Remove the duplicate bold "Play Biergarten Fest Free..." paragraph that
#    used to sit right before the final italic paragraph.
# ---------------------------------------------------------------------------
$boldDupPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$boldDupPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Replace the text of the final italic paragraph with the new image
#    prompt, again via a direct Range so the leading empty run + <w:i/>
#    formatting on the run survive untouched.
# ---------------------------------------------------------------------------
$italicPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$oldItalicText = "Discover Biergarten Fest slot game, inspired by Oktoberfest. Enjoy Wilds, Scatters, and bonus round. Play for free, and experience immersive graphics."
$italicRunStart = $italicPara.Range.Start
$italicRunEnd = $italicRunStart + $oldItalicText.Length
$italicRun = $d.Range($italicRunStart, $italicRunEnd)

$newItalicText = 'Prompt: Create a cartoon-style feature image for "Biergarten Fest" that features a happy Maya warrior with glasses. The image should be in bright, vibrant colors and should incorporate elements of German beer festivals such as beer mugs, pretzels, and traditional German attire. The Maya warrior should be depicted holding a beer mug while enjoying the festivities. The background should feature a lively beer garden with market stalls, traditional fairs, and people having a good time. Make sure to include the Biergarten Fest logo and convey the fun and cheerful atmosphere of the game through the image.'
$italicRun.Text = $newItalicText

Write-Host "Done"
